$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D10", "D16", "D19", "D20", "D23", "D27", "D29", "D31", "D32", "D34", "D36", "D38", "D40", "D42", "D47", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "71.234.30"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "2.624.08"
$ws.Range("E3").Value = "  +4.65%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "606.35"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "181.00"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "2.622.88"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +15.19%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E15").Value = "  +8.72%  "
$ws.Range("D16").Value = "26.58"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "71.319.98"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").Value = "2.629.00"
$ws.Range("E18").Value = "  +6.11%  "
$ws.Range("D19").Value = "383.03"
$ws.Range("E19").Value = "  +9.13%  "
$ws.Range("D20").Value = "7.90"
$ws.Range("E20").Value = "  +6.31%  "
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "72.21"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("E24").Value = "  +6.04%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +9.32%  "
$ws.Range("D27").Value = "9.64"
$ws.Range("E27").Value = "  +5.79%  "
$ws.Range("D28").Value = "2.759.51"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +7.33%  "
$ws.Range("D31").Value = "544.28"
$ws.Range("E31").Value = "  +6.69%  "
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("E33").Value = "  +5.33%  "
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "165.62"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "19.19"
$ws.Range("E39").Value = "  +8.06%  "
$ws.Range("D40").Value = "18.98"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  +5.00%  "
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +9.49%  "
$ws.Range("E43").Value = "  +5.47%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("D47").Value = "154.55"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("E49").Value = "  +6.10%  "
$ws.Range("D50").Value = "0.532"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").Value = "0.0₆0265"
$ws.Range("E51").Value = "  +2.37%  "
